$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.5
$ws.Range("I3").Value = 6
$ws.Range("K3").Value = 15
$ws.Range("N3").Value = 1.53
$ws.Range("O3").Value = 2.5
$ws.Range("AH3").Value = 67

# Row 5
$ws.Range("N5").Value = 2.3
$ws.Range("O5").Value = 1.6

# Row 8
$ws.Range("J8").Value = 1.05
$ws.Range("K8").Value = 11
$ws.Range("N8").Value = 1.95
$ws.Range("O8").Value = 1.85

# Row 11
$ws.Range("G11").Value = 1.65
$ws.Range("H11").Value = 3.5
$ws.Range("I11").Value = 5.1
$ws.Range("L11").Value = 1.42
$ws.Range("M11").Value = 2.47
$ws.Range("N11").Value = 2.22
$ws.Range("O11").Value = 1.52
$ws.Range("P11").Value = 1.5
$ws.Range("Q11").Value = 2.27
$ws.Range("R11").Value = 2.18
$ws.Range("S11").Value = 1.53
$ws.Range("T11").Value = 5.1
$ws.Range("U11").Value = 6.4
$ws.Range("V11").Value = 8.75
$ws.Range("W11").Value = 11.5
$ws.Range("X11").Value = 16
$ws.Range("Y11").Value = 40
$ws.Range("Z11").Value = 7.4
$ws.Range("AA11").Value = 7.1
$ws.Range("AB11").Value = 23
$ws.Range("AC11").Value = 150
$ws.Range("AE11").Value = 10.5
$ws.Range("AF11").Value = 27
$ws.Range("AG11").Value = 18
$ws.Range("AH11").Value = 100
$ws.Range("AI11").Value = 70
$ws.Range("AJ11").Value = 90

# Row 12
$ws.Range("G12").Value = 2.3
$ws.Range("H12").Value = 2.95
$ws.Range("I12").Value = 3.15
$ws.Range("L12").Value = 1.47
$ws.Range("M12").Value = 2.35
$ws.Range("N12").Value = 2.32
$ws.Range("O12").Value = 1.47
$ws.Range("P12").Value = 1.55
$ws.Range("Q12").Value = 2.15
$ws.Range("R12").Value = 1.98
$ws.Range("S12").Value = 1.65
$ws.Range("T12").Value = 6
$ws.Range("U12").Value = 9.75
$ws.Range("V12").Value = 9.75
$ws.Range("W12").Value = 23
$ws.Range("X12").Value = 23
$ws.Range("Y12").Value = 40
$ws.Range("Z12").Value = 6.7
$ws.Range("AA12").Value = 5.9
$ws.Range("AB12").Value = 17.5
$ws.Range("AC12").Value = 110
$ws.Range("AE12").Value = 7.5
$ws.Range("AF12").Value = 15
$ws.Range("AG12").Value = 11.75
$ws.Range("AH12").Value = 45
$ws.Range("AI12").Value = 35
$ws.Range("AJ12").Value = 50

# Row 13
$ws.Range("N13").Value = 1.57
$ws.Range("O13").Value = 2.35

# Row 15
$ws.Range("G15").Value = 1.52
$ws.Range("H15").Value = 3.7
$ws.Range("I15").Value = 5.6
$ws.Range("N15").Value = 1.82
$ws.Range("P15").Value = 1.37
$ws.Range("Q15").Value = 2.5
$ws.Range("T15").Value = 5.4
$ws.Range("U15").Value = 5.9
$ws.Range("W15").Value = 8.75
$ws.Range("X15").Value = 10.5
$ws.Range("Z15").Value = 9.75
$ws.Range("AA15").Value = 6.4
$ws.Range("AE15").Value = 11.75
$ws.Range("AF15").Value = 27
$ws.Range("AG15").Value = 15
$ws.Range("AH15").Value = 80
$ws.Range("AI15").Value = 45

# Row 16
$ws.Range("G16").Value = 1.85
$ws.Range("H16").Value = 3.15
$ws.Range("I16").Value = 4
$ws.Range("N16").Value = 2.02
$ws.Range("O16").Value = 1.62
$ws.Range("P16").Value = 1.38
$ws.Range("Q16").Value = 2.45
$ws.Range("T16").Value = 5.5
$ws.Range("U16").Value = 7.1
$ws.Range("V16").Value = 7.1
$ws.Range("W16").Value = 13
$ws.Range("X16").Value = 12.5
$ws.Range("Y16").Value = 23
$ws.Range("Z16").Value = 7.9
$ws.Range("AA16").Value = 5.4
$ws.Range("AB16").Value = 13
$ws.Range("AC16").Value = 60
$ws.Range("AD16").Value = 450
$ws.Range("AE16").Value = 8.25
$ws.Range("AF16").Value = 17
$ws.Range("AG16").Value = 11.25
$ws.Range("AH16").Value = 50
$ws.Range("AI16").Value = 32
$ws.Range("AJ16").Value = 37

# Row 17
$ws.Range("G17").Value = 1.75
$ws.Range("H17").Value = 3.25
$ws.Range("I17").Value = 4.3
$ws.Range("N17").Value = 1.93
$ws.Range("O17").Value = 1.7
$ws.Range("P17").Value = 1.37
$ws.Range("Q17").Value = 2.5
$ws.Range("U17").Value = 6.7
$ws.Range("V17").Value = 6.9
$ws.Range("W17").Value = 11.75
$ws.Range("X17").Value = 12
$ws.Range("Z17").Value = 8.75
$ws.Range("AA17").Value = 5.6
$ws.Range("AD17").Value = 350
$ws.Range("AE17").Value = 10
$ws.Range("AF17").Value = 20
$ws.Range("AG17").Value = 11.5
$ws.Range("AH17").Value = 55
$ws.Range("AI17").Value = 32
$ws.Range("AJ17").Value = 35

# Row 20
$ws.Range("G20").Value = 2.4
$ws.Range("I20").Value = 2.9
$ws.Range("V20").Value = 9.5
$ws.Range("X20").Value = 19
$ws.Range("AA20").Value = 6.5
$ws.Range("AJ20").Value = 34

# Row 21
$ws.Range("G21").Value = 3.1
$ws.Range("I21").Value = 2.45
$ws.Range("T21").Value = 8.5
$ws.Range("U21").Value = 15
$ws.Range("V21").Value = 12
$ws.Range("W21").Value = 34
$ws.Range("X21").Value = 29
$ws.Range("AE21").Value = 7
$ws.Range("AF21").Value = 11
$ws.Range("AH21").Value = 23
$ws.Range("AI21").Value = 21

# Row 24
$ws.Range("H24").Value = 6
$ws.Range("I24").Value = 9.5
$ws.Range("J24").Value = 1.02
$ws.Range("K24").Value = 12
$ws.Range("N24").Value = 1.4
$ws.Range("O24").Value = 2.75
$ws.Range("U24").Value = 7.5
$ws.Range("V24").Value = 9.5
$ws.Range("AA24").Value = 12

# Row 26
$ws.Range("G26").Value = 1.95
$ws.Range("H26").Value = 3.9
$ws.Range("I26").Value = 3.1
$ws.Range("J26").Value = 1.02
$ws.Range("K26").Value = 12
$ws.Range("T26").Value = 13
$ws.Range("V26").Value = 9.5
$ws.Range("Y26").Value = 17
$ws.Range("AE26").Value = 19
$ws.Range("AG26").Value = 13
$ws.Range("AJ26").Value = 23

# Row 30
$ws.Range("H30").Value = 5.5
$ws.Range("K30").Value = 9.75
$ws.Range("L30").Value = 1.13
$ws.Range("M30").Value = 5.1
$ws.Range("N30").Value = 1.42
$ws.Range("O30").Value = 2.65
$ws.Range("P30").Value = 1.24
$ws.Range("Q30").Value = 3.6
$ws.Range("R30").Value = 1.82
$ws.Range("S30").Value = 1.88
$ws.Range("T30").Value = 9.25
$ws.Range("U30").Value = 7.3
$ws.Range("V30").Value = 9
$ws.Range("W30").Value = 8.25
$ws.Range("X30").Value = 10
$ws.Range("Y30").Value = 23
$ws.Range("Z30").Value = 9.75
$ws.Range("AA30").Value = 11.5
$ws.Range("AE30").Value = 35

# Row 31
$ws.Range("G31").Value = 3.95
$ws.Range("H31").Value = 3.9
$ws.Range("I31").Value = 1.75
$ws.Range("P31").Value = 1.34
$ws.Range("Q31").Value = 3
$ws.Range("T31").Value = 13
$ws.Range("U31").Value = 23
$ws.Range("Y31").Value = 37
$ws.Range("AA31").Value = 7.6
$ws.Range("AB31").Value = 14.5
$ws.Range("AC31").Value = 60
$ws.Range("AF31").Value = 9
$ws.Range("AH31").Value = 14
$ws.Range("AI31").Value = 13
